# Natmi following Dr Hou advice
# Re-derive the Fgf17 -> Fgfr2 LR-pair table across the full
# Sending-cluster x Target-cluster grid (FAPs/sCs x ECs/FAPs/M2/sCs),
# replacing the 2-row placeholder with the complete 8-row result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("FAPs", "Fgf17", "Fgfr2", "ECs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7530926191740721, 0.7530926191740722, 2, 0.6666666666666666, 0.7050016666666666, 2.115005, 0.2162043499394578, 0.2162043499394578, 1.008673849566111, 9.078064646095, 0.1628219001727339, 0.162821900172734),
    @("FAPs", "Fgf17", "Fgfr2", "FAPs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7530926191740721, 0.7530926191740722, 3, 1, 2.337487333333333, 7.012461999999999, 0.716842176819984, 0.716842176819984, 3.344335848130888, 30.09902263317799, 0.5398485524758051, 0.5398485524758051),
    @("FAPs", "Fgf17", "Fgfr2", "M2", 3, 1, 1.430739666666667, 4.292218999999999, 0.7530926191740721, 0.7530926191740722, 1, 0.3333333333333333, 0.007121333333333334, 0.021364, 0.002183914332167809, 0.002183914332167809, 0.01018877407955555, 0.09169896671599999, 0.00164468976446405, 0.00164468976446405),
    @("FAPs", "Fgf17", "Fgfr2", "sCs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7530926191740721, 0.7530926191740722, 3, 1, 0.2112013333333334, 0.6336040000000001, 0.0647695589083904, 0.06476955890839041, 0.3021741252528889, 2.719567127276, 0.04877747676106908, 0.0487774767610691),
    @("sCs", "Fgf17", "Fgfr2", "ECs", 2, 0.6666666666666666, 0.4690793333333333, 1.407238, 0.2469073808259278, 0.2469073808259279, 2, 0.6666666666666666, 0.7050016666666666, 2.115005, 0.2162043499394578, 0.2162043499394578, 0.3307017117988889, 2.97631540619, 0.05338244976672387, 0.05338244976672389),
    @("sCs", "Fgf17", "Fgfr2", "FAPs", 2, 0.6666666666666666, 0.4690793333333333, 1.407238, 0.2469073808259278, 0.2469073808259279, 3, 1, 2.337487333333333, 7.012461999999999, 0.716842176819984, 0.716842176819984, 1.096466999995111, 9.868202999955999, 0.1769936243441789, 0.1769936243441789),
    @("sCs", "Fgf17", "Fgfr2", "M2", 2, 0.6666666666666666, 0.4690793333333333, 1.407238, 0.2469073808259278, 0.2469073808259279, 1, 0.3333333333333333, 0.007121333333333334, 0.021364, 0.002183914332167809, 0.002183914332167809, 0.003340470292444445, 0.030064232632, 0.000539224567703759, 0.0005392245677037592),
    @("sCs", "Fgf17", "Fgfr2", "sCs", 2, 0.6666666666666666, 0.4690793333333333, 1.407238, 0.2469073808259278, 0.2469073808259279, 3, 1, 0.2112013333333334, 0.6336040000000001, 0.0647695589083904, 0.06476955890839041, 0.09907018063911112, 0.8916316257520001, 0.01599208214732131, 0.01599208214732132),
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rows[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $rowValues[$c]
    }
}
